$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P: "Final Project Presentation" (header in row 2, next to Midterm 2 in O2)
$ws.Range("P2").Value = "Final Project Presentation"

# Widen the new column P to roughly match the authored width (25.11 characters).
# Note: the COM column-width model here quantizes to whole pixels, so we pick
# the input that lands on the closest representable stored width.
$ws.Columns.Item(16).ColumnWidth = 24.28

# Mark the Final Project Presentation as completed ("1") for the three
# students who submitted it.
$ws.Range("P11").Formula = "=1"
$ws.Range("P13").Formula = "=1"
$ws.Range("P18").Formula = "=1"

# Update the active selection to match the authored state.
$ws.Range("P11").Select() | Out-Null
